$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("URLsWithDomains")

# Remove the "about.html", "our-team.html", "careers.html" rows (A2:A4).
# The remaining rows shift up, and "about-solar.html" becomes the new A2.
$ws.Range("A2:A4").EntireRow.Delete()

# Update the (now stale) hidden AutoFilter defined name range so it keeps
# the same padding below the data as before (it shrinks by the 3 deleted rows).
foreach ($n in $wb.Names) {
    if ($n.Name -eq "URLsWithDomains!_FilterDatabase") {
        $n.RefersTo = "=URLsWithDomains!`$A`$1:`$A`$27"
    }
}

# Reset the view: select E4 (instead of D4) and let the top-left scroll
# position return to the default (A1).
$ws.Range("E4").Select()
